# Update the cryptos price/volume snapshot (GitHub Actions refresh).
# Columns: D = Price (text, dotted thousands/decimals as scraped), E = Volume(1h) (text, padded "  +/-x.xx%  ").
# Price values that look like plain numbers are entered with a leading "'" so Excel
# keeps them as text (matching the original inlineStr/string cell type) instead of
# coercing them into numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.792.62"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "2.353.30"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'240.03"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("E6").Value = "  -1.87%  "

$ws.Range("D7").Value = "'73.24"
$ws.Range("E7").Value = "  -1.60%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "'0.601"

$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("D11").Value = "'60.77"
$ws.Range("E11").Value = "  +6.14%  "

$ws.Range("D12").Value = "'35.04"
$ws.Range("E12").Value = "  +8.41%  "

$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("D14").Value = "'7.15"
$ws.Range("E14").Value = "  -2.44%  "

$ws.Range("E15").Value = "  -2.68%  "

$ws.Range("D16").Value = "'0.908"
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("D17").Value = "2.365.08"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").Value = "43.753.52"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("E21").Value = "  -3.32%  "

$ws.Range("D22").Value = "'252.68"
$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").Value = "  +3.04%  "

$ws.Range("D25").Value = "'1.88"
$ws.Range("E25").Value = "  -5.09%  "

$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("D27").Value = "'10.47"
$ws.Range("E27").Value = "  -3.20%  "

$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("D29").Value = "'175.52"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").Value = "'22.24"
$ws.Range("E30").Value = "  -2.35%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  -2.54%  "

$ws.Range("D33").Value = "'0.0746"
$ws.Range("E33").Value = "  -3.81%  "

$ws.Range("E34").Value = "  -4.22%  "

$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("E37").Value = "  +4.06%  "

$ws.Range("E38").Value = "  +1.75%  "

$ws.Range("D39").Value = "'0.0277"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").Value = "'5.53"
$ws.Range("E40").Value = "  +16.16%  "

$ws.Range("D41").Value = "'64.93"
$ws.Range("E41").Value = "  +11.69%  "

$ws.Range("D42").Value = "'19.83"
$ws.Range("E42").Value = "  +3.63%  "

$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("E44").Value = "  -6.22%  "

$ws.Range("D45").Value = "'9.02"
$ws.Range("E45").Value = "  -1.16%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("E47").Value = "  -3.24%  "

$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("D50").Value = "'97.67"
$ws.Range("E50").Value = "  -2.73%  "
